# repull data, push all data, mean calculation
# Update the "dSF" column (F) with repulled values for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0
$ws.Range("F3").Value = -3
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = -1
$ws.Range("F6").Value = -1
$ws.Range("F7").Value = 7
$ws.Range("F8").Value = -6
$ws.Range("F10").Value = -5
$ws.Range("F14").Value = -3
$ws.Range("F17").Value = 1
$ws.Range("F19").Value = 3
